# Change default delimiter character from '~' to '.' in the example workbook.
# Strings containing the old '~' delimiter are rewritten in place, touching the
# "nasa~a_low~#" / "nasa~a_high~#" columns (H:U) first, then the
# "elements~H" / "elements~O" columns (C:D) -- matching the order the values
# were originally edited in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Update-Delimiter($cell) {
    $v = $cell.Value2
    if ($v -is [string] -and $v.Contains("~")) {
        $cell.Value2 = $v.Replace("~", ".")
    }
}

# nasa~a_low~0..6 and nasa~a_high~0..6 live in columns H through U of row 1
for ($c = 8; $c -le 21; $c++) {
    Update-Delimiter $ws.Cells.Item(1, $c)
}

# elements~H and elements~O live in columns C and D of row 1
for ($c = 3; $c -le 4; $c++) {
    Update-Delimiter $ws.Cells.Item(1, $c)
}

# Update the saved view state: scroll the frozen pane back to column B,
# and leave the active selection on E1.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 2
[void]$ws.Range("E1").Select()
